$wb = $excel.ActiveWorkbook

# Update the title on the "Inputs and Outputs" sheet to "Results Summary and Inputs"
$ws = $wb.Worksheets.Item("Inputs and Outputs")
$ws.Range("A1").Value = "Results Summary and Inputs"

# Make this sheet the active/selected tab, and move the selection to A2
$ws.Activate()
$ws.Range("A2").Select()
